# Add model_no_12 and model_no_13 rows to the model_parameters sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 12) down onto the
# two new rows (13 and 14) so the new cells pick up the same cell styles
# (s="5" for A/C/D/E, s="1" for B) that the rest of the table uses.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 13: model_no_12 ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 8994859
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = "'0.0242"
$ws.Range("E13").Value = 140

# --- Row 14: model_no_13 ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 2580523
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = "'0.0236"
$ws.Range("E14").Value = 164

# Re-apply row 12's formatting onto the "Best val_loss" cells so that typing
# the value as text (to match the existing column's text-stored numbers)
# doesn't leave the cells on a different style than the rest of the column.
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Match the final selection left by the author in the saved workbook.
$ws.Range("C14").Select()
